$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.579.91"
$ws.Range("E2").Value = "  -0.78%  "
$ws.Range("D3").Value = "1.877.63"
$ws.Range("E3").Value = "  -1.18%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.66"
$ws.Range("E5").Value = "  -4.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4869"
$ws.Range("E7").Value = "  -2.10%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2907"
$ws.Range("E8").Value = "  -2.28%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06660"
$ws.Range("E9").Value = "  -2.52%  "
$ws.Range("D10").Value = "1.879.00"
$ws.Range("E10").Value = "  -1.08%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "16.59"
$ws.Range("E11").Value = "  -3.72%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07236"
$ws.Range("E12").Value = "  -1.26%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "88.91"
$ws.Range("E13").Value = "  -2.94%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.992"
$ws.Range("E14").Value = "  -2.03%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6519"
$ws.Range("E15").Value = "  -4.02%  "
$ws.Range("D16").Value = "30.518.58"
$ws.Range("E16").Value = "  -0.96%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000007855"
$ws.Range("E17").Value = "  -2.44%  "
$ws.Range("E18").Value = "  +0.20%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.94"
$ws.Range("E19").Value = "  -3.56%  "
$ws.Range("D20").Value = "2.120.44"
$ws.Range("E20").Value = "  -1.48%  "
$ws.Range("E21").Value = "  -0.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.705"
$ws.Range("E22").Value = "  -3.56%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "191.75"
$ws.Range("E23").Value = "  +6.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.095"
$ws.Range("E24").Value = "  -0.11%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.265"
$ws.Range("E25").Value = "  -0.97%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "157.75"
$ws.Range("E26").Value = "  +1.57%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.32"
$ws.Range("E27").Value = "  -2.25%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.823"
$ws.Range("E28").Value = "  -6.31%  "
$ws.Range("E29").Value = "  +0.74%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.237"
$ws.Range("E30").Value = "  -2.96%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08996"
$ws.Range("E31").Value = "  +0.38%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.918"
$ws.Range("E32").Value = "  -3.38%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05129"
$ws.Range("E33").Value = "  -3.16%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7216"
$ws.Range("E34").Value = "  -4.19%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.077"
$ws.Range("E35").Value = "  -5.93%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.691"
$ws.Range("E36").Value = "  -0.25%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01813"
$ws.Range("E37").Value = "  -5.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.659"
$ws.Range("E38").Value = "  -2.32%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9170"
$ws.Range("E39").Value = "  -2.22%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.049"
$ws.Range("E40").Value = "  -6.12%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4372"
$ws.Range("E41").Value = "  -0.32%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "104.64"
$ws.Range("E42").Value = "  -1.27%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9952"
$ws.Range("E43").Value = "  -0.65%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.699"
$ws.Range("E44").Value = "  -2.81%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1329"
$ws.Range("E45").Value = "  -3.50%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.340"
$ws.Range("E46").Value = "  -5.23%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4010"
$ws.Range("E47").Value = "  +2.33%  "
$ws.Range("E48").Value = "  -0.43%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.603"
$ws.Range("E49").Value = "  +0.24%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.402"
$ws.Range("E50").Value = "  +0.61%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "33.06"
$ws.Range("E51").Value = "  -1.62%  "
